# Stundenerfassung.xlsx - add two new time-tracking entries and switch the
# active worksheet back to "Stundenerfassung" (from "Wochen").

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Stundenerfassung")

# --- Row 108: 14.08.2017 (42961), ETIC2 / Design View Model, 5h ---
# Insert a copy of the last existing row (107) so the new row inherits the
# exact same cell formatting (date format, etc.) instead of getting a blank
# default style, then overwrite the values.
$ws1.Rows("107:107").Copy() | Out-Null
$ws1.Rows("108:108").Insert(-4121) | Out-Null
$ws1.Range("A108").Value = 42961
$ws1.Range("B108").Value = "ETIC2"
$ws1.Range("C108").Value = "Design View Model"
$ws1.Range("D108").Value = 5

# --- Row 109: 15.08.2017 (42962), ETIC2 / Design View Model, 4h ---
$ws1.Rows("108:108").Copy() | Out-Null
$ws1.Rows("109:109").Insert(-4121) | Out-Null
$ws1.Range("A109").Value = 42962
$ws1.Range("B109").Value = "ETIC2"
$ws1.Range("C109").Value = "Design View Model"
$ws1.Range("D109").Value = 4

# Switch the active tab back to "Stundenerfassung" and leave the selection on
# the last touched cell.
$ws1.Activate() | Out-Null
$ws1.Range("E106").Select() | Out-Null
